$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.103.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.884.59"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4974"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2918"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06630"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "1.879.89"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07196"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6664"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.851"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "30.092.95"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007826"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9974"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "2.121.12"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9970"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.768"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.618"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.163"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.913"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.379"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08680"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.950"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7053"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.655"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.202"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9356"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01650"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4199"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.514"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1262"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05712"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.348"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.99%  "
